# Apply updated cryptocurrency price/volume data to Sheet1
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range('D2')
$c.NumberFormat = '@'
$c.Value = '36.500.76'
$c.Style = 'Normal'

$c = $ws.Range('D3')
$c.NumberFormat = '@'
$c.Value = '1.985.67'
$c.Style = 'Normal'

$c = $ws.Range('E3')
$c.NumberFormat = '@'
$c.Value = '  -3.33%  '
$c.Style = 'Normal'

$c = $ws.Range('E4')
$c.NumberFormat = '@'
$c.Value = '  -0.08%  '
$c.Style = 'Normal'

$c = $ws.Range('D5')
$c.NumberFormat = '@'
$c.Value = '244.68'
$c.Style = 'Normal'

$c = $ws.Range('E5')
$c.NumberFormat = '@'
$c.Value = '  -3.09%  '
$c.Style = 'Normal'

$c = $ws.Range('E6')
$c.NumberFormat = '@'
$c.Value = '  -3.84%  '
$c.Style = 'Normal'

$c = $ws.Range('D7')
$c.NumberFormat = '@'
$c.Value = '59.54'
$c.Style = 'Normal'

$c = $ws.Range('E7')
$c.NumberFormat = '@'
$c.Value = '  -9.10%  '
$c.Style = 'Normal'

$c = $ws.Range('E8')
$c.NumberFormat = '@'
$c.Value = '  +0.00%  '
$c.Style = 'Normal'

$c = $ws.Range('E9')
$c.NumberFormat = '@'
$c.Value = '  -2.10%  '
$c.Style = 'Normal'

$c = $ws.Range('D10')
$c.NumberFormat = '@'
$c.Value = '57.65'
$c.Style = 'Normal'

$c = $ws.Range('E10')
$c.NumberFormat = '@'
$c.Value = '  -3.54%  '
$c.Style = 'Normal'

$c = $ws.Range('D11')
$c.NumberFormat = '@'
$c.Value = '0.0823'
$c.Style = 'Normal'

$c = $ws.Range('E11')
$c.NumberFormat = '@'
$c.Value = '  +6.83%  '
$c.Style = 'Normal'

$c = $ws.Range('D13')
$c.NumberFormat = '@'
$c.Value = '23.90'
$c.Style = 'Normal'

$c = $ws.Range('E13')
$c.NumberFormat = '@'
$c.Value = '  +5.65%  '
$c.Style = 'Normal'

$c = $ws.Range('E14')
$c.NumberFormat = '@'
$c.Value = '  -5.92%  '
$c.Style = 'Normal'

$c = $ws.Range('E15')
$c.NumberFormat = '@'
$c.Value = '  -5.40%  '
$c.Style = 'Normal'

$c = $ws.Range('D16')
$c.NumberFormat = '@'
$c.Value = '2.275.63'
$c.Style = 'Normal'

$c = $ws.Range('E16')
$c.NumberFormat = '@'
$c.Value = '  -3.44%  '
$c.Style = 'Normal'

$c = $ws.Range('D17')
$c.NumberFormat = '@'
$c.Value = '5.49'
$c.Style = 'Normal'

$c = $ws.Range('E17')
$c.NumberFormat = '@'
$c.Value = '  -2.01%  '
$c.Style = 'Normal'

$c = $ws.Range('D18')
$c.NumberFormat = '@'
$c.Value = '1.984.03'
$c.Style = 'Normal'

$c = $ws.Range('E18')
$c.NumberFormat = '@'
$c.Value = '  -3.48%  '
$c.Style = 'Normal'

$c = $ws.Range('D19')
$c.NumberFormat = '@'
$c.Value = '36.367.48'
$c.Style = 'Normal'

$c = $ws.Range('E19')
$c.NumberFormat = '@'
$c.Value = '  -2.67%  '
$c.Style = 'Normal'

$c = $ws.Range('D20')
$c.NumberFormat = '@'
$c.Value = '70.31'
$c.Style = 'Normal'

$c = $ws.Range('E20')
$c.NumberFormat = '@'
$c.Value = '  -4.60%  '
$c.Style = 'Normal'

$c = $ws.Range('D21')
$c.NumberFormat = '@'
$c.Value = '0.0₃0864'
$c.Style = 'Normal'

$c = $ws.Range('E21')
$c.NumberFormat = '@'
$c.Value = '  -1.47%  '
$c.Style = 'Normal'

$c = $ws.Range('D22')
$c.NumberFormat = '@'
$c.Value = '5.33'
$c.Style = 'Normal'

$c = $ws.Range('E22')
$c.NumberFormat = '@'
$c.Value = '  -2.75%  '
$c.Style = 'Normal'

$c = $ws.Range('D23')
$c.NumberFormat = '@'
$c.Value = '234.26'
$c.Style = 'Normal'

$c = $ws.Range('E23')
$c.NumberFormat = '@'
$c.Value = '  -2.55%  '
$c.Style = 'Normal'

$c = $ws.Range('E24')
$c.NumberFormat = '@'
$c.Value = '  -0.02%  '
$c.Style = 'Normal'

$c = $ws.Range('D25')
$c.NumberFormat = '@'
$c.Value = '2.61'
$c.Style = 'Normal'

$c = $ws.Range('E25')
$c.NumberFormat = '@'
$c.Value = '  -0.68%  '
$c.Style = 'Normal'

$c = $ws.Range('D26')
$c.NumberFormat = '@'
$c.Value = '2.31'
$c.Style = 'Normal'

$c = $ws.Range('E26')
$c.NumberFormat = '@'
$c.Value = '  -3.99%  '
$c.Style = 'Normal'

$c = $ws.Range('D27')
$c.NumberFormat = '@'
$c.Value = '10.06'
$c.Style = 'Normal'

$c = $ws.Range('E27')
$c.NumberFormat = '@'
$c.Value = '  -0.89%  '
$c.Style = 'Normal'

$c = $ws.Range('D28')
$c.NumberFormat = '@'
$c.Value = '162.33'
$c.Style = 'Normal'

$c = $ws.Range('E28')
$c.NumberFormat = '@'
$c.Value = '  +0.05%  '
$c.Style = 'Normal'

$c = $ws.Range('D29')
$c.NumberFormat = '@'
$c.Value = '19.86'
$c.Style = 'Normal'

$c = $ws.Range('E29')
$c.NumberFormat = '@'
$c.Value = '  -0.76%  '
$c.Style = 'Normal'

$c = $ws.Range('D30')
$c.NumberFormat = '@'
$c.Value = '0.132'
$c.Style = 'Normal'

$c = $ws.Range('E30')
$c.NumberFormat = '@'
$c.Value = '  +11.32%  '
$c.Style = 'Normal'

$c = $ws.Range('E31')
$c.NumberFormat = '@'
$c.Value = '  -1.80%  '
$c.Style = 'Normal'

$c = $ws.Range('E32')
$c.NumberFormat = '@'
$c.Value = '  -1.26%  '
$c.Style = 'Normal'

$c = $ws.Range('E33')
$c.NumberFormat = '@'
$c.Value = '  -6.93%  '
$c.Style = 'Normal'

$c = $ws.Range('D34')
$c.NumberFormat = '@'
$c.Value = '0.0631'
$c.Style = 'Normal'

$c = $ws.Range('E34')
$c.NumberFormat = '@'
$c.Value = '  +0.58%  '
$c.Style = 'Normal'

$c = $ws.Range('E35')
$c.NumberFormat = '@'
$c.Value = '  -6.12%  '
$c.Style = 'Normal'

$c = $ws.Range('D36')
$c.NumberFormat = '@'
$c.Value = '6.32'
$c.Style = 'Normal'

$c = $ws.Range('E36')
$c.NumberFormat = '@'
$c.Value = '  +4.53%  '
$c.Style = 'Normal'

$c = $ws.Range('E37')
$c.NumberFormat = '@'
$c.Value = '  +0.00%  '
$c.Style = 'Normal'

$c = $ws.Range('E38')
$c.NumberFormat = '@'
$c.Value = '  -7.19%  '
$c.Style = 'Normal'

$c = $ws.Range('D39')
$c.NumberFormat = '@'
$c.Value = '1.78'
$c.Style = 'Normal'

$c = $ws.Range('E39')
$c.NumberFormat = '@'
$c.Value = '  -2.87%  '
$c.Style = 'Normal'

$c = $ws.Range('D40')
$c.NumberFormat = '@'
$c.Value = '3.05'
$c.Style = 'Normal'

$c = $ws.Range('E40')
$c.NumberFormat = '@'
$c.Value = '  +2.59%  '
$c.Style = 'Normal'

$c = $ws.Range('E41')
$c.NumberFormat = '@'
$c.Value = '  -0.75%  '
$c.Style = 'Normal'

$c = $ws.Range('D42')
$c.NumberFormat = '@'
$c.Value = '0.0964'
$c.Style = 'Normal'

$c = $ws.Range('E42')
$c.NumberFormat = '@'
$c.Value = '  -7.74%  '
$c.Style = 'Normal'

$c = $ws.Range('E43')
$c.NumberFormat = '@'
$c.Value = '  -3.71%  '
$c.Style = 'Normal'

$c = $ws.Range('E44')
$c.NumberFormat = '@'
$c.Value = '  -2.40%  '
$c.Style = 'Normal'

$c = $ws.Range('E45')
$c.NumberFormat = '@'
$c.Value = '  -4.99%  '
$c.Style = 'Normal'

$c = $ws.Range('D46')
$c.NumberFormat = '@'
$c.Value = '92.86'
$c.Style = 'Normal'

$c = $ws.Range('E46')
$c.NumberFormat = '@'
$c.Value = '  -4.15%  '
$c.Style = 'Normal'

$c = $ws.Range('D47')
$c.NumberFormat = '@'
$c.Value = '16.24'
$c.Style = 'Normal'

$c = $ws.Range('E47')
$c.NumberFormat = '@'
$c.Value = '  -5.94%  '
$c.Style = 'Normal'

$c = $ws.Range('D48')
$c.NumberFormat = '@'
$c.Value = '1.380.38'
$c.Style = 'Normal'

$c = $ws.Range('E48')
$c.NumberFormat = '@'
$c.Value = '  -2.90%  '
$c.Style = 'Normal'

$c = $ws.Range('D49')
$c.NumberFormat = '@'
$c.Value = '7.52'
$c.Style = 'Normal'

$c = $ws.Range('E49')
$c.NumberFormat = '@'
$c.Value = '  -5.81%  '
$c.Style = 'Normal'

$c = $ws.Range('E50')
$c.NumberFormat = '@'
$c.Value = '  -3.28%  '
$c.Style = 'Normal'

$c = $ws.Range('D51')
$c.NumberFormat = '@'
$c.Value = '45.43'
$c.Style = 'Normal'

$c = $ws.Range('E51')
$c.NumberFormat = '@'
$c.Value = '  -2.84%  '
$c.Style = 'Normal'

